# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 10:05"

# Row 5 - Rusia
$ws.Range("B5").Value = 290678
$ws.Range("C5").Value = 8926
$ws.Range("D5").Value = 70209
$ws.Range("E5").Value = 217747
$ws.Range("G5").Value = 91
$ws.Range("H5").Value = 2722

# Row 88 - Estonia
$ws.Range("B88").Value = 1784
$ws.Range("C88").Value = 10
$ws.Range("E88").Value = 782
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 64

# Row 92 - Eslovaquia
$ws.Range("B92").Value = 1495
$ws.Range("C92").Value = 1
$ws.Range("D92").Value = 1185
$ws.Range("E92").Value = 282

# Row 106 - Sri Lanka
$ws.Range("D106").Value = 559
$ws.Range("E106").Value = 413
